$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 668.5641000000001
$ws.Range("I15").Value = 668.5641000000001
$ws.Range("K15").Value = 2005.6923
$ws.Range("M15").Value = -1836.6923

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 3506.6667
$ws.Range("I40").Value = 4492.8125
$ws.Range("J40").Value = 2072.2727
$ws.Range("K40").Value = 4492.8125
$ws.Range("L40").Value = 2072.2727
$ws.Range("M40").Value = -4317.8125
$ws.Range("N40").Value = -2422.2727

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 2100
$ws.Range("I51").Value = 1800
$ws.Range("J51").Value = 2400
$ws.Range("K51").Value = 1800
$ws.Range("L51").Value = 2400
$ws.Range("M51").Value = -1316
$ws.Range("N51").Value = -3368

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 7163.1763
$ws.Range("I64").Value = 3799.1667
$ws.Range("J64").Value = 8998.091
$ws.Range("K64").Value = 3799.1667
$ws.Range("L64").Value = 8998.091
$ws.Range("M64").Value = -3551.1667
$ws.Range("N64").Value = -9494.091

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 7163.1763
$ws.Range("I67").Value = 3799.1667
$ws.Range("J67").Value = 8998.091
$ws.Range("K67").Value = 3799.1667
$ws.Range("L67").Value = 8998.091
$ws.Range("M67").Value = -2941.1667
$ws.Range("N67").Value = -10714.091

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 6719.6
$ws.Range("I74").Value = 4585.2856
$ws.Range("J74").Value = 9436
$ws.Range("K74").Value = 4585.2856
$ws.Range("L74").Value = 9436
$ws.Range("M74").Value = -3649.2856
$ws.Range("N74").Value = -11308

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 6719.6
$ws.Range("I77").Value = 4585.2856
$ws.Range("J77").Value = 9436
$ws.Range("K77").Value = 22926.428
$ws.Range("L77").Value = 47180
$ws.Range("M77").Value = -18246.428
$ws.Range("N77").Value = -56540

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 6097.5
$ws.Range("J100").Value = 8160.2
$ws.Range("L100").Value = 8160.2
$ws.Range("N100").Value = -9242.200000000001

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 3537.0176
$ws.Range("I132").Value = 2604.327
$ws.Range("K132").Value = 7812.981000000001
$ws.Range("M132").Value = -5282.981000000001

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 55557476
$ws.Range("I137").Value = 100000570
$ws.Range("J137").Value = 3610
$ws.Range("K137").Value = 300001710
$ws.Range("L137").Value = 10830
$ws.Range("M137").Value = -299999160
$ws.Range("N137").Value = -15930

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2507.5
$ws.Range("I138").Value = 1300
$ws.Range("J138").Value = 3436.3462
$ws.Range("K138").Value = 3900
$ws.Range("L138").Value = 10309.0386
$ws.Range("M138").Value = 1240
$ws.Range("N138").Value = -20589.0386

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3772.4285
$ws.Range("I61").Value = 2892.182
$ws.Range("K61").Value = 2892.182
$ws.Range("M61").Value = -2680.182

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 3149.818
$ws.Range("I102").Value = 3106.25
$ws.Range("K102").Value = 3106.25
$ws.Range("M102").Value = -1484.25

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1939.26
$ws.Range("I132").Value = 1486.5333
$ws.Range("K132").Value = 4459.5999
$ws.Range("M132").Value = -1929.5999

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3772.4285
$ws.Range("I136").Value = 2892.182
$ws.Range("K136").Value = 8676.545999999998
$ws.Range("M136").Value = -6126.545999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1884.6666
$ws.Range("I134").Value = 1149.4468
$ws.Range("K134").Value = 3448.3404
$ws.Range("M134").Value = -913.3404

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3646.8096
$ws.Range("I31").Value = 1325.091
$ws.Range("J31").Value = 6200.7
$ws.Range("K31").Value = 1325.091
$ws.Range("L31").Value = 6200.7
$ws.Range("M31").Value = -1030.091
$ws.Range("N31").Value = -6790.7

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3646.8096
$ws.Range("I34").Value = 1325.091
$ws.Range("J34").Value = 6200.7
$ws.Range("K34").Value = 1325.091
$ws.Range("L34").Value = 6200.7
$ws.Range("M34").Value = -1123.091
$ws.Range("N34").Value = -6604.7

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 2127.8076
$ws.Range("I134").Value = 1927.6522
$ws.Range("K134").Value = 5782.9566
$ws.Range("M134").Value = -3247.9566

$ws = $wb.Worksheets.Item("CUL")
# Row 9: Jack of All Plates / Jack-o'-lantern
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 10563.154
$ws.Range("I68").Value = 2875
$ws.Range("J68").Value = 13980.111
$ws.Range("K68").Value = 8625
$ws.Range("L68").Value = 41940.333
$ws.Range("M68").Value = -7814
$ws.Range("N68").Value = -43562.333

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 10563.154
$ws.Range("I71").Value = 2875
$ws.Range("J71").Value = 13980.111
$ws.Range("K71").Value = 25875
$ws.Range("L71").Value = 125820.999
$ws.Range("M71").Value = -21819
$ws.Range("N71").Value = -133932.999

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 4470
$ws.Range("I80").Value = 4065
$ws.Range("J80").Value = 4773.75
$ws.Range("K80").Value = 12195
$ws.Range("L80").Value = 14321.25
$ws.Range("M80").Value = -11259
$ws.Range("N80").Value = -16193.25

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 4470
$ws.Range("I83").Value = 4065
$ws.Range("J83").Value = 4773.75
$ws.Range("K83").Value = 36585
$ws.Range("L83").Value = 42963.75
$ws.Range("M83").Value = -31905
$ws.Range("N83").Value = -52323.75

# Row 116: On a Full Stomach / Sausage Links
$ws.Range("H116").Value = 475
$ws.Range("I116").Value = 475
$ws.Range("K116").Value = 1425
$ws.Range("M116").Value = 2017

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 34511.15
$ws.Range("I122").Value = 42878.92
$ws.Range("J122").Value = 11267.333
$ws.Range("K122").Value = 128636.76
$ws.Range("L122").Value = 33801.999
$ws.Range("M122").Value = -126186.76
$ws.Range("N122").Value = -38701.999

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 5649.8
$ws.Range("I126").Value = 2812.5
$ws.Range("J126").Value = 16999
$ws.Range("K126").Value = 8437.5
$ws.Range("L126").Value = 50997
$ws.Range("M126").Value = -5967.5
$ws.Range("N126").Value = -55937

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 14297129
$ws.Range("I132").Value = 18529944
$ws.Range("K132").Value = 55589832
$ws.Range("M132").Value = -55587302

$ws = $wb.Worksheets.Item("LTW")
# Row 11: A Thorn in One's Hide / Leather Mitts
$ws.Range("H11").Value = 2007
$ws.Range("J11").Value = 2007
$ws.Range("L11").Value = 2007
$ws.Range("N11").Value = -2287

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2449.739
$ws.Range("I132").Value = 2325.9524
$ws.Range("K132").Value = 6977.8572
$ws.Range("M132").Value = -4447.8572

# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 58165
$ws.Range("J133").Value = 58165
$ws.Range("L133").Value = 58165
$ws.Range("N133").Value = -63225

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4055.0645
$ws.Range("I136").Value = 2229.5454
$ws.Range("K136").Value = 6688.6362
$ws.Range("M136").Value = -4138.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 20: Read the Fine Print / Cotton Shepherd's Tunic
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 4351202.5
$ws.Range("I132").Value = 5130597.5
$ws.Range("K132").Value = 15391792.5
$ws.Range("M132").Value = -15389262.5

# Row 139: Cruel Climates / Rroneek Serge Trousers of Gathering
$ws.Range("H139").Value = 86100
$ws.Range("I139").Value = 87750
$ws.Range("J139").Value = 85000
$ws.Range("K139").Value = 87750
$ws.Range("L139").Value = 87750
$ws.Range("M139").Value = -82610
$ws.Range("N139").Value = -95280
